$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B29").Value = "ESEHTWO-941"
$ws.Rows.Item(29).RowHeight = 90
Write-Host ("height=" + $ws.Rows.Item(29).RowHeight)
